$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.8100456314743383
$ws.Cells.Item(2, 3).Value = 0.2139355548678168
$ws.Cells.Item(2, 5).Value = 0.1180885000064151
$ws.Cells.Item(2, 6).Value = 0.4443680307746121
$ws.Cells.Item(2, 7).Value = 0.002380715531325821
$ws.Cells.Item(2, 9).Value = 0.282251155278157
$ws.Cells.Item(2, 13).Value = 0.3415018856951235
$ws.Cells.Item(2, 15).Value = 1.237202940415017

$ws.Cells.Item(3, 2).Value = 0.7079256111462655
$ws.Cells.Item(3, 3).Value = 0.1889873452860229
$ws.Cells.Item(3, 5).Value = 0.1134358236564665
$ws.Cells.Item(3, 6).Value = 0.3878228170618172
$ws.Cells.Item(3, 7).Value = 0.002383361759907004
$ws.Cells.Item(3, 9).Value = 0.2904908242850794
$ws.Cells.Item(3, 13).Value = 0.3029369837983751
$ws.Cells.Item(3, 15).Value = 1.256615068974668

$ws.Cells.Item(4, 2).Value = 0.6450042764207069
$ws.Cells.Item(4, 3).Value = 0.1735931867904981
$ws.Cells.Item(4, 5).Value = 0.1107153921229767
$ws.Cells.Item(4, 6).Value = 0.3531389305169483
$ws.Cells.Item(4, 7).Value = 0.00238507146575101
$ws.Cells.Item(4, 9).Value = 0.2958889709655157
$ws.Cells.Item(4, 13).Value = 0.2792820313245201
$ws.Cells.Item(4, 15).Value = 1.270069308549935

$ws.Cells.Item(5, 2).Value = 0.6193097748508478
$ws.Cells.Item(5, 3).Value = 0.1673012155555966
$ws.Cells.Item(5, 5).Value = 0.1096407594327715
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.002385789601615425
$ws.Cells.Item(5, 9).Value = 0.2981737049354241
$ws.Cells.Item(5, 13).Value = 0.2696486179753208
$ws.Cells.Item(5, 15).Value = 1.275936454532555

$ws.Cells.Item(6, 2).Value = 0.6150400352093754
$ws.Cells.Item(6, 3).Value = 0.1662553178174164
$ws.Cells.Item(6, 5).Value = 0.1094643596679283
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.002385910142949602
$ws.Cells.Item(6, 9).Value = 0.2985582053158176
$ws.Cells.Item(6, 13).Value = 0.2680493735698661
$ws.Cells.Item(6, 15).Value = 1.276933862092889

$ws.Cells.Item(7, 2).Value = 0.6446579660883458
$ws.Cells.Item(7, 3).Value = 0.173508406423025
$ws.Cells.Item(7, 5).Value = 0.1107007621445248
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.002385081064079973
$ws.Cells.Item(7, 9).Value = 0.2959194401291541
$ws.Cells.Item(7, 13).Value = 0.2791520865844532
$ws.Cells.Item(7, 15).Value = 1.270146880477583

$ws.Cells.Item(8, 2).Value = 0.7748810616644164
$ws.Cells.Item(8, 3).Value = 0.2053493566115776
$ws.Cells.Item(8, 5).Value = 0.116455769791564
$ws.Cells.Item(8, 6).Value = 0.4248636149813478
$ws.Cells.Item(8, 7).Value = 0.002381610366809352
$ws.Cells.Item(8, 9).Value = 0.2850216706646904
$ws.Cells.Item(8, 13).Value = 0.3281997760878212
$ws.Cells.Item(8, 15).Value = 1.243576678678536

$ws.Cells.Item(9, 2).Value = 1.028451974918596
$ws.Cells.Item(9, 3).Value = 0.2671752339249451
$ws.Cells.Item(9, 5).Value = 0.1288380069051271
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.002375475102967428
$ws.Cells.Item(9, 9).Value = 0.2663532087977849
$ws.Cells.Item(9, 13).Value = 0.4245726733431212
$ws.Cells.Item(9, 15).Value = 1.20372435167404

$ws.Cells.Item(10, 2).Value = 1.213598073173443
$ws.Cells.Item(10, 3).Value = 0.3122123748455294
$ws.Cells.Item(10, 5).Value = 0.1386255262341933
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.002371372260207526
$ws.Cells.Item(10, 9).Value = 0.2543014884139243
$ws.Cells.Item(10, 13).Value = 0.4955011030209562
$ws.Cells.Item(10, 15).Value = 1.18201243765391

$ws.Cells.Item(11, 2).Value = 1.297564813878751
$ws.Cells.Item(11, 3).Value = 0.3326148552239943
$ws.Cells.Item(11, 5).Value = 0.1432327538483378
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.002369592764406604
$ws.Cells.Item(11, 9).Value = 0.2491837447280449
$ws.Cells.Item(11, 13).Value = 0.5277971025173542
$ws.Cells.Item(11, 15).Value = 1.173799517412888

$ws.Cells.Item(12, 2).Value = 1.329322571574721
$ws.Cells.Item(12, 3).Value = 0.3403282240847716
$ws.Cells.Item(12, 5).Value = 0.1450000193504053
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.002368931345714801
$ws.Cells.Item(12, 9).Value = 0.2472985309681306
$ws.Cells.Item(12, 13).Value = 0.5400311233909889
$ws.Cells.Item(12, 15).Value = 1.170930490999552

$ws.Cells.Item(13, 2).Value = 1.322484714104235
$ws.Cells.Item(13, 3).Value = 0.3386675785760076
$ws.Cells.Item(13, 5).Value = 0.1446183961151988
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.002369073241919747
$ws.Cells.Item(13, 9).Value = 0.2477021935150248
$ws.Cells.Item(13, 13).Value = 0.5373961221418853
$ws.Cells.Item(13, 15).Value = 1.171537639448005

$ws.Cells.Item(14, 2).Value = 1.300178327181584
$ws.Cells.Item(14, 3).Value = 0.333249693470492
$ws.Cells.Item(14, 5).Value = 0.1433776927698602
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.002369538100146131
$ws.Cells.Item(14, 9).Value = 0.2490275876883423
$ws.Cells.Item(14, 13).Value = 0.5288035187973037
$ws.Cells.Item(14, 15).Value = 1.173558640271921

$ws.Cells.Item(15, 2).Value = 1.286509942470218
$ws.Cells.Item(15, 3).Value = 0.3299294297726476
$ws.Cells.Item(15, 5).Value = 0.1426206806278358
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.002369824457302225
$ws.Cells.Item(15, 9).Value = 0.249846311182969
$ws.Cells.Item(15, 13).Value = 0.5235408478016126
$ws.Cells.Item(15, 15).Value = 1.174828002682005

$ws.Cells.Item(16, 2).Value = 1.208105333489073
$ws.Cells.Item(16, 3).Value = 0.3108772802558804
$ws.Cells.Item(16, 5).Value = 0.1383275780838176
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.002371490298397592
$ws.Cells.Item(16, 9).Value = 0.2546433107268076
$ws.Cells.Item(16, 13).Value = 0.4933910819755027
$ws.Cells.Item(16, 15).Value = 1.182582804175965

$ws.Cells.Item(17, 2).Value = 1.159939605668569
$ws.Cells.Item(17, 3).Value = 0.2991673362128608
$ws.Cells.Item(17, 5).Value = 0.1357337998102963
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.002372534454422944
$ws.Cells.Item(17, 9).Value = 0.2576797430662356
$ws.Cells.Item(17, 13).Value = 0.4749028437014431
$ws.Cells.Item(17, 15).Value = 1.187767548663899

$ws.Cells.Item(18, 2).Value = 1.132211846227904
$ws.Cells.Item(18, 3).Value = 0.2924240937570062
$ws.Cells.Item(18, 5).Value = 0.1342564804776814
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.00237314320903102
$ws.Cells.Item(18, 9).Value = 0.2594605210075756
$ws.Cells.Item(18, 13).Value = 0.4642717576357285
$ws.Cells.Item(18, 15).Value = 1.190906189539959

$ws.Cells.Item(19, 2).Value = 1.122819613453714
$ws.Cells.Item(19, 3).Value = 0.2901395862096194
$ws.Cells.Item(19, 5).Value = 0.1337587746416418
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.002373350730271507
$ws.Cells.Item(19, 9).Value = 0.2600693432527041
$ws.Cells.Item(19, 13).Value = 0.4606727497487384
$ws.Cells.Item(19, 15).Value = 1.1919957071173

$ws.Cells.Item(20, 2).Value = 1.165069436568842
$ws.Cells.Item(20, 3).Value = 0.3004147095823839
$ws.Cells.Item(20, 5).Value = 0.1360084032388897
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.002372422455593836
$ws.Cells.Item(20, 9).Value = 0.2573529571067237
$ws.Cells.Item(20, 13).Value = 0.4768706532334193
$ws.Cells.Item(20, 15).Value = 1.187199412521622

$ws.Cells.Item(21, 2).Value = 1.306731315032664
$ws.Cells.Item(21, 3).Value = 0.3348414027081219
$ws.Cells.Item(21, 5).Value = 0.1437415010925207
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.002369401223004484
$ws.Cells.Item(21, 9).Value = 0.2486368526487563
$ws.Cells.Item(21, 13).Value = 0.5313272609905226
$ws.Cells.Item(21, 15).Value = 1.172958467924587

$ws.Cells.Item(22, 2).Value = 1.399089656457647
$ws.Cells.Item(22, 3).Value = 0.3572675836865642
$ws.Cells.Item(22, 5).Value = 0.148927466551946
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.002367499141820099
$ws.Cells.Item(22, 9).Value = 0.2432480803641361
$ws.Cells.Item(22, 13).Value = 0.5669423818644788
$ws.Cells.Item(22, 15).Value = 1.165056897076425

$ws.Cells.Item(23, 2).Value = 1.349817479332955
$ws.Cells.Item(23, 3).Value = 0.3453051673787968
$ws.Cells.Item(23, 5).Value = 0.1461474331197223
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.002368507707039308
$ws.Cells.Item(23, 9).Value = 0.2460959056647578
$ws.Cells.Item(23, 13).Value = 0.5479317144278895
$ws.Cells.Item(23, 15).Value = 1.169144912261459

$ws.Cells.Item(24, 2).Value = 1.162750354029356
$ws.Cells.Item(24, 3).Value = 0.299850806443942
$ws.Cells.Item(24, 5).Value = 0.1358842118270687
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.002372473063977705
$ws.Cells.Item(24, 9).Value = 0.2575005877512648
$ws.Cells.Item(24, 13).Value = 0.4759810126192718
$ws.Cells.Item(24, 15).Value = 1.187455775268361

$ws.Cells.Item(25, 2).Value = 0.9600526307479527
$ws.Cells.Item(25, 3).Value = 0.2505166946359623
$ws.Cells.Item(25, 5).Value = 0.1253686339503446
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.002377063478087167
$ws.Cells.Item(25, 9).Value = 0.271112518101738
$ws.Cells.Item(25, 13).Value = 0.3984801369161488
$ws.Cells.Item(25, 15).Value = 1.213183624856427
